# Update "想去人数" (interested-count) figures to match the regenerated
# gh-pages data snapshot (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F3").Value = 2852   # was 2834
$wsExpo.Range("F5").Value = 11     # was 8

# Sheet "演出" (performances)
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F2").Value = 117    # was 116

# Sheet "全部类型" (all types, combined)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 117     # was 116
$wsAll.Range("F7").Value = 2852    # was 2834
$wsAll.Range("F10").Value = 11     # was 8
